# Fix casing of Building Material names (BOM) used by the wastage-factor
# lookup table so they match the naming convention used elsewhere
# (sentence case instead of Title Case), per "fix: all bom and wastage
# factor fixes for compatibility".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value  = "Normalweight concrete, 4000 psi"
$ws.Range("A8").Value  = "Normalweight concrete, 6000 psi"
$ws.Range("A9").Value  = "Normalweight concrete, 5000 psi"
$ws.Range("A12").Value = "Steel curtain wall system"
$ws.Range("A14").Value = "5`" Mineral wool insulation"
$ws.Range("A18").Value = "Aluminum curtain wall system"
$ws.Range("A20").Value = "Type N mortar"
$ws.Range("A22").Value = "4`" Mineral wool insulation"
$ws.Range("A30").Value = "3/8`" Synthetic stucco"
$ws.Range("A32").Value = "Glass fiber reinforced concrete (GFRC) Panel"
$ws.Range("A34").Value = "4.5`" Mineral wool insulation"
$ws.Range("A36").Value = "Galvanized steel backer tray"
$ws.Range("A37").Value = "Steel, sheet"
$ws.Range("A38").Value = "Thermal break"
$ws.Range("A39").Value = "3.5`" Mineral wool insulation"
$ws.Range("A40").Value = "6`" x 1' Tulipwood/Poplar lumber"
$ws.Range("A43").Value = "Formed steel sheet"
$ws.Range("A44").Value = "Enamel paint"
$ws.Range("A45").Value = "Stainless steel fasteners"
$ws.Range("A46").Value = "Galvanized steel support"

# Leave the scroll position where the author's session landed: cell A33
# selected (matches the saved view in the updated workbook).
$ws.Range("A33").Select()
